# Creacion de policies a traves de Data Driven
# Adds a new "Policy" worksheet (Data Driven template) with 10 sample policy rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Restore / update the selection on the existing AAPolicy sheet
# ---------------------------------------------------------------------------
$wsPolicyOld = $wb.Worksheets.Item("AAPolicy")
$wsPolicyOld.Activate()
$wsPolicyOld.Range("E1:H1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1

# ---------------------------------------------------------------------------
# 2) Add the new "Policy" worksheet after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Policy"

# Copy the header look (bold font + yellow fill) used on the other sheets
$wsPolicyOld.Range("A1").Copy() | Out-Null
$ws.Range("A1:O1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Policy_Name"
$ws.Range("B1").Value = "Policy_Currency"
$ws.Range("C1").Value = "From_Date"
$ws.Range("D1").Value = "To_Date"
$ws.Range("E1").Value = "Covars_Number"
$ws.Range("F1").Value = "Covar1"
$ws.Range("G1").Value = "Covar2"
$ws.Range("H1").Value = "Covar3"
$ws.Range("I1").Value = "Covar4"
$ws.Range("J1").Value = "Covar5"
$ws.Range("K1").Value = "Covar6"
$ws.Range("L1").Value = "Covar7"
$ws.Range("M1").Value = "Covar8"
$ws.Range("N1").Value = "Covar9"
$ws.Range("O1").Value = "Covar10"

# ---------------------------------------------------------------------------
# 4) Template / first data row (row 2) - all the "shared" values
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "United States Dollar"

$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = "1/1/2018"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "12/31/2018"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Net Pay"
$ws.Range("G2").Value = "zzzAlimony Received"
$ws.Range("H2").Value = "zzzAnnual Paid Premium"
$ws.Range("I2").Value = "zzzAnnuities"
$ws.Range("J2").Value = "zzzAuto Loan Advance"
$ws.Range("K2").Value = "zzzAutomobile Expense"
$ws.Range("L2").Value = "zzzBase Salary"
$ws.Range("M2").Value = "zzzBonus - Current Year"
$ws.Range("N2").Value = "zzzEducation Trip"
$ws.Range("O2").Value = "zzzEntertainment Allowance"

# Propagate the template row's formatting (esp. the C2 date format) down
# to rows 3-11 by copying the formats only, so no duplicate styles are
# generated.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 5) Replicate the template values (columns B, C, D, E-O) for rows 3-11
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("B$r").Value = "United States Dollar"
    $ws.Range("C$r").Value = "1/1/2018"

    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = "12/31/2018"
    $ws.Range("D$r").Style = "Normal"

    $ws.Range("E$r").Value = 10
    $ws.Range("F$r").Value = "Net Pay"
    $ws.Range("G$r").Value = "zzzAlimony Received"
    $ws.Range("H$r").Value = "zzzAnnual Paid Premium"
    $ws.Range("I$r").Value = "zzzAnnuities"
    $ws.Range("J$r").Value = "zzzAuto Loan Advance"
    $ws.Range("K$r").Value = "zzzAutomobile Expense"
    $ws.Range("L$r").Value = "zzzBase Salary"
    $ws.Range("M$r").Value = "zzzBonus - Current Year"
    $ws.Range("N$r").Value = "zzzEducation Trip"
    $ws.Range("O$r").Value = "zzzEntertainment Allowance"
}

# ---------------------------------------------------------------------------
# 6) Fill in the policy names (column A) for every data row, last
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Automated Policy28"
$ws.Range("A3").Value = "Automated Policy29"
$ws.Range("A4").Value = "Automated Policy30"
$ws.Range("A5").Value = "Automated Policy31"
$ws.Range("A6").Value = "Automated Policy32"
$ws.Range("A7").Value = "Automated Policy33"
$ws.Range("A8").Value = "Automated Policy34"
$ws.Range("A9").Value = "Automated Policy35"
$ws.Range("A10").Value = "Automated Policy36"
$ws.Range("A11").Value = "Automated Policy37"

# ---------------------------------------------------------------------------
# 7) Column widths (best fit approximation)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 18.0
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334
$ws.Columns.Item(5).ColumnWidth = 14.5
$ws.Columns.Item(6).ColumnWidth = 7.0
$ws.Columns.Item(7).ColumnWidth = 19.0
$ws.Columns.Item(8).ColumnWidth = 22.333333333333332
$ws.Columns.Item(9).ColumnWidth = 11.333333333333334
$ws.Columns.Item(10).ColumnWidth = 19.833333333333332
$ws.Columns.Item(11).ColumnWidth = 21.666666666666668
$ws.Columns.Item(12).ColumnWidth = 12.666666666666666
$ws.Columns.Item(13).ColumnWidth = 21.166666666666668
$ws.Columns.Item(14).ColumnWidth = 15.5
$ws.Columns.Item(15).ColumnWidth = 25.833333333333332

# ---------------------------------------------------------------------------
# 8) Final selection / active cell for the new sheet
# ---------------------------------------------------------------------------
$ws.Range("A10").Select() | Out-Null
